# Updated cryptos list: refresh Price (D) and Volume(1h) (E) columns for rows 2-51.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Each entry: row, new Price text (or $null if unchanged), new Volume(1h) text (or $null if unchanged).
$updates = @(
    ,@(2, '60.328.83', '  -4.10%  ')
    ,@(3, '2.977.14', '  -5.91%  ')
    ,@(4, $null, '  -0.06%  ')
    ,@(5, '566.68', '  -3.92%  ')
    ,@(6, '124.17', '  -7.28%  ')
    ,@(7, $null, '  +0.13%  ')
    ,@(8, '2.971.67', '  -6.03%  ')
    ,@(9, '0.497', '  -3.45%  ')
    ,@(10, $null, '  -6.23%  ')
    ,@(11, '5.12', '  -2.33%  ')
    ,@(12, '0.435', '  -4.06%  ')
    ,@(13, $null, '  -5.85%  ')
    ,@(14, '32.43', '  -6.99%  ')
    ,@(15, $null, '  -0.63%  ')
    ,@(16, '3.468.27', '  -5.89%  ')
    ,@(17, '60.350.58', '  -4.04%  ')
    ,@(18, '2.971.54', '  -6.25%  ')
    ,@(19, '6.13', '  -6.37%  ')
    ,@(20, '424.75', '  -7.75%  ')
    ,@(21, $null, '  -6.28%  ')
    ,@(22, '0.658', '  -5.39%  ')
    ,@(23, '7.07', '  -7.33%  ')
    ,@(24, '12.80', '  -4.36%  ')
    ,@(25, '78.51', '  -5.36%  ')
    ,@(26, $null, '  +0.12%  ')
    ,@(27, $null, '  -0.25%  ')
    ,@(28, $null, '  -6.21%  ')
    ,@(29, $null, '  -7.56%  ')
    ,@(30, $null, '  -7.79%  ')
    ,@(31, '25.08', '  -7.37%  ')
    ,@(32, $null, '  -11.30%  ')
    ,@(33, $null, '  -9.74%  ')
    ,@(34, $null, '  -5.02%  ')
    ,@(35, '0.942', '  -8.79%  ')
    ,@(36, $null, '  -4.65%  ')
    ,@(37, '49.31', '  -3.85%  ')
    ,@(38, '0.0₃0647', '  -7.61%  ')
    ,@(39, $null, '  -8.22%  ')
    ,@(40, '7.84', '  -3.05%  ')
    ,@(41, '0.109', '  -3.09%  ')
    ,@(42, '376.08', '  -6.34%  ')
    ,@(43, '2.629.16', '  -5.72%  ')
    ,@(44, $null, '  -8.80%  ')
    ,@(45, $null, '  +0.03%  ')
    ,@(46, $null, '  -6.62%  ')
    ,@(47, '118.88', '  -4.71%  ')
    ,@(48, $null, '  -7.48%  ')
    ,@(49, '0.106', '  -4.79%  ')
    ,@(50, '23.21', '  -7.90%  ')
    ,@(51, '31.14', '  -9.11%  ')
)

foreach ($u in $updates) {
    $row = $u[0]
    $priceText = $u[1]
    $volumeText = $u[2]

    if ($null -ne $priceText) {
        $priceCell = $ws.Cells.Item($row, 4)
        # Force text storage (matches the sheet's existing text-typed Price column) even
        # for values that parse as plain decimals, so trailing zeros survive (e.g. '12.80').
        $needsTextForce = $priceText -match '^-?[0-9]+(\.[0-9]+)?$'
        if ($needsTextForce) {
            $priceCell.Value = "'" + $priceText
            $priceCell.Style = "Normal"
        } else {
            $priceCell.Value = $priceText
        }
    }

    if ($null -ne $volumeText) {
        $ws.Cells.Item($row, 5).Value = $volumeText
    }
}
